# Wireframes document: bump "Version 1." -> "Version 2."
# The target OOXML splits the word "Version" into two runs ("Versi" + "on"),
# changes the " 1." run to " 2" (dropping the trailing period), and appends
# a brand-new run containing "." after the existing _GoBack bookmark.

$d = $word.ActiveDocument

# --- Step 1: split the "Version" run into "Versi" | "on" -------------------
# Temporarily dropping a bookmark at the split point forces the engine to
# break the run in two without introducing any explicit run formatting
# (rPr), matching the diff's plain <w:r><w:t>Versi</w:t></w:r> style split.
$splitPoint = $d.Range(5, 5)
$d.Bookmarks.Add("TempRunSplit", $splitPoint)
$d.Bookmarks("TempRunSplit").Delete()

# --- Step 2: turn " 1." into " 2" -------------------------------------------
# "Version" occupies characters 0-7, so the following run " 1." starts at 7;
# its digit+period ("1.") live at 8-10. Replace that with just "2".
$versionNumber = $d.Range(8, 10)
$versionNumber.Text = "2"

# --- Step 3: re-add the trailing period as its own run after the bookmark --
$tailInsertionPoint = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$tailInsertionPoint.InsertAfter(".")
